$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.825.10'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.758.07'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5072'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.21%  '
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2628'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06194'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("D11").Value = '1.752.46'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06929'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.23%  '
$ws.Range("E13").Value = '  +7.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6030'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.451'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '25.861.99'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006820'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.59%  '
$ws.Range("E21").Value = '  +3.64%  '
$ws.Range("D22").Value = '1.977.21'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +5.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.154'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.78%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.458'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.63%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.811'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08221'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.693'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.390'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04362'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +1.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.002'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6000'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.735'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("E40").Value = '  +4.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.931'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3814'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7446'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.868'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05486'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.25%  '
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.943'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("E51").Value = '  -0.05%  '
